# Updated cryptos list values (prices + % change) per upstream diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.967.82"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").Value = "2.305.24"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "'304.10"
$ws.Range("E5").Value = "  +1.16%  "
$ws.Range("D6").Value = "'97.09"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.502"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.500"
$ws.Range("D10").Value = "'35.33"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("D11").Value = "'0.0788"
$ws.Range("E11").Value = "  +0.01%  "
$ws.Range("D12").Value = "'18.80"
$ws.Range("E12").Value = "  +4.95%  "
$ws.Range("E13").Value = "  +1.53%  "
$ws.Range("D14").Value = "'6.89"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "2.665.13"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "2.299.06"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("D17").Value = "'0.782"
$ws.Range("D18").Value = "42.861.37"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "'12.59"
$ws.Range("E19").Value = "  -1.53%  "
$ws.Range("E20").Value = "  -0.87%  "
$ws.Range("D21").Value = "'6.04"
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("D22").Value = "'67.67"
$ws.Range("E22").Value = "  -0.31%  "
$ws.Range("D23").Value = "'236.90"
$ws.Range("E23").Value = "  -1.67%  "
$ws.Range("D24").Value = "'2.17"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D26").Value = "'2.42"
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "'24.88"
$ws.Range("E27").Value = "  -2.04%  "
$ws.Range("D28").Value = "'2.37"
$ws.Range("E28").Value = "  +17.17%  "
$ws.Range("D29").Value = "'165.80"
$ws.Range("E29").Value = "  +0.10%  "
$ws.Range("D30").Value = "'9.05"
$ws.Range("E30").Value = "  +0.11%  "
$ws.Range("D31").Value = "'32.79"
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'18.14"
$ws.Range("E33").Value = "  +6.29%  "
$ws.Range("E34").Value = "  -0.92%  "
$ws.Range("E35").Value = "  -8.08%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  +0.48%  "
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "'1.75"
$ws.Range("E40").Value = "  +1.48%  "
$ws.Range("E41").Value = "  -0.67%  "
$ws.Range("D42").Value = "1.996.93"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "'10.30"
$ws.Range("E44").Value = "  +1.36%  "
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "'17.89"
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("D47").Value = "'2.76"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "2.531.17"
$ws.Range("E48").Value = "  +0.48%  "
$ws.Range("D49").Value = "'53.38"
$ws.Range("E49").Value = "  -0.27%  "
$ws.Range("E50").Value = "  -2.41%  "
$ws.Range("D51").Value = "'71.75"
$ws.Range("E51").Value = "  -0.39%  "
